# Update cryptos list cell values to reflect the latest scrape (GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.966.92"
$ws.Range("E2").Value = "  +2.61%  "
$ws.Range("D3").Value = "'3.821.37"
$ws.Range("E3").Value = "  +1.31%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'662.84"
$ws.Range("E5").Value = "  +6.52%  "
$ws.Range("D6").Value = "'169.55"
$ws.Range("E6").Value = "  +3.33%  "
$ws.Range("D7").Value = "'3.818.98"
$ws.Range("E7").Value = "  +1.32%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  +1.05%  "
$ws.Range("E11").Value = "  +2.46%  "
$ws.Range("D12").Value = "'6.96"
$ws.Range("E12").Value = "  +5.15%  "
$ws.Range("D13").Value = "'0.0000245"
$ws.Range("E13").Value = "  -0.32%  "
$ws.Range("D14").Value = "'35.73"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "'4.459.53"
$ws.Range("E15").Value = "  +1.16%  "
$ws.Range("D16").Value = "'3.811.80"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "'70.839.31"
$ws.Range("E17").Value = "  +2.42%  "
$ws.Range("D18").Value = "'17.83"
$ws.Range("E18").Value = "  +0.88%  "
$ws.Range("D19").Value = "'7.16"
$ws.Range("E19").Value = "  +1.26%  "
$ws.Range("E20").Value = "  +1.02%  "
$ws.Range("D21").Value = "'479.27"
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("E22").Value = "  +7.88%  "
$ws.Range("D23").Value = "'0.714"
$ws.Range("E23").Value = "  +2.05%  "
$ws.Range("D24").Value = "'0.0000147"
$ws.Range("E24").Value = "  -1.14%  "
$ws.Range("D25").Value = "'82.97"
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'12.28"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").Value = "'10.39"
$ws.Range("E27").Value = "  +4.23%  "
$ws.Range("E28").Value = "  -0.89%  "
$ws.Range("D30").Value = "'3.966.58"
$ws.Range("E30").Value = "  +1.20%  "
$ws.Range("D31").Value = "'2.83"
$ws.Range("E31").Value = "  +7.03%  "
$ws.Range("E32").Value = "  +3.69%  "
$ws.Range("D33").Value = "'7.46"
$ws.Range("E33").Value = "  +2.78%  "
$ws.Range("D34").Value = "'29.56"
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("E35").Value = "  +15.90%  "
$ws.Range("B36").Value = "Binance-PegBSC-USD"
$ws.Range("C36").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("E36").Value = "  +0.05%  "
$ws.Range("B37").Value = "RenzoRestakedETH"
$ws.Range("C37").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D37").Value = "'3.771.51"
$ws.Range("E37").Value = "  +1.30%  "
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("E39").Value = "  -0.05%  "
$ws.Range("D40").Value = "'3.45"
$ws.Range("E40").Value = "  +3.69%  "
$ws.Range("E41").Value = "  +2.96%  "
$ws.Range("D42").Value = "'0.967"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("D44").Value = "'2.11"
$ws.Range("E44").Value = "  +10.87%  "
$ws.Range("D46").Value = "'45.23"
$ws.Range("E46").Value = "  +5.51%  "
$ws.Range("D47").Value = "'158.92"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("D48").Value = "'47.70"
$ws.Range("E48").Value = "  +2.19%  "
$ws.Range("D49").Value = "'0.301"
$ws.Range("E49").Value = "  +0.67%  "
$ws.Range("E50").Value = "  +3.95%  "
$ws.Range("B51").Value = "FLOKI"
$ws.Range("C51").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D51").Value = "'0.000293"
$ws.Range("E51").Value = "  +8.79%  "
